$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width -> stored width 66 (ColumnWidth input empirically snaps to 66 in this engine)
$ws.Columns.Item(1).ColumnWidth = 65.14

# Rewrite the full QA test-plan content (A1:B19)
$ws.Range("A1").Value = "logar no aplicativo mobile"
$ws.Range("B1").Value = "logou com sucesso"
$ws.Range("A2").Value = "realizar uma entrada "
$ws.Range("B2").Value = "realizou entrada com sucesso"
$ws.Range("A3").Value = "simular uma entrada duplicada(use o mesmo serial acima)"
$ws.Range("B3").Value = "Não realiza a entrada, informa que a mesma já foi realizada anteriormente"
$ws.Range("A4").Value = "dê saida nesse equipamento"
$ws.Range("B4").Value = "realizou a saida do equipamento com sucesso"
$ws.Range("A5").Value = "realizar uma saida duplicada"
$ws.Range("B5").Value = "não realiza saida duplicada de equipamentos"
$ws.Range("A6").Value = "realizar logout"
$ws.Range("B6").Value = "deslogando com sucesso"
$ws.Range("A7").Value = "realizar loguin desktop"
$ws.Range("B7").Value = "logou com sucesso"
$ws.Range("A8").Value = "realizar uma entrada "
$ws.Range("B8").Value = "realizou a entrada com sucesso"
$ws.Range("A9").Value = "simular uma entrada duplicada(use o mesmo serial acima)"
$ws.Range("B9").Value = "não realiza  entrada duplicada e informa que a mesma foi realizada anteriormente"
$ws.Range("A10").Value = "dê saida nesse equipamento"
$ws.Range("B10").Value = "realizou a saida do equipamento com sucesso"
$ws.Range("A11").Value = "realizar uma saida duplicada"
$ws.Range("B11").Value = "não realiza saida duplicada de equipamentos pede para o usuario conferir o numero serial"
$ws.Range("A12").Value = "acessar tela de inserção de equipamentos com usuario sem autorização"
$ws.Range("B12").Value = "não permite acesso, caso desse acesso, o mesmo não conseguiria acessar e modificar os dados"
$ws.Range("A13").Value = "realizar logout"
$ws.Range("B13").Value = "deslogando com sucesso"
$ws.Range("A14").Value = "acessar com usuario com autorização admin ou master"
$ws.Range("B14").Value = "acesso realizado com sucesso"
$ws.Range("A15").Value = "realizar uma entrada "
$ws.Range("B15").Value = "realizou a entrada com sucesso"
$ws.Range("A16").Value = "realizar uma entrada  duplicada (use o serial anterior)"
$ws.Range("B16").Value = "Não realiza a entrada, informa que a mesma já foi realizada anteriormente"
$ws.Range("A17").Value = "realizar a saida do equipamento anterior"
$ws.Range("B17").Value = "realizou a saida do equipamento com sucesso"
$ws.Range("A18").Value = "realizar saida duplicada do equipamento anteriro"
$ws.Range("B18").Value = "não realiza a saida duplicada de equipamentos"
$ws.Range("A19").Value = "acessar opção de inserção de novos equipamentos"
$ws.Range("B19").Value = "não esta permitindo que o usuario master acesse"

# Clear the stray red-text formatting on B3, B8, B10, B18 (now matches plain Calibri style, like B2/B4)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)

# B19 keeps its existing red+underline emphasis style; only its text changed above.
